# Added a "4wk low sales check" to the forecast summary workbook.
# This updates the Inventory Coverage (H) and Seasonality Index (L) values
# on the "Forecast Comparison" sheet (now accounting for the trailing 4-week
# low-sales look-back), a couple of MyForecast (D) quantities that shifted
# as a result, and the two dependent roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# Row 2 (W10)
$ws1.Range("H2").Value = 29
$ws1.Range("L2").Value = 1.17

# Row 3 (W11)
$ws1.Range("D3").Value = 1
$ws1.Range("H3").Value = 28
$ws1.Range("L3").Value = 1.19

# Row 4 (W12)
$ws1.Range("H4").Value = 13.5
$ws1.Range("L4").Value = 1.16

# Row 5 (W13)
$ws1.Range("H5").Value = 12.5
$ws1.Range("L5").Value = 0.96

# Row 6 (W14)
$ws1.Range("H6").Value = 11.5
$ws1.Range("L6").Value = 0.85

# Row 7 (W15)
$ws1.Range("D7").Value = 2
$ws1.Range("H7").Value = 10.5
$ws1.Range("L7").Value = 0.91

# Row 8 (W16)
$ws1.Range("D8").Value = 2
$ws1.Range("H8").Value = 9.5
$ws1.Range("L8").Value = 0.83

# Row 9 (W17)
$ws1.Range("D9").Value = 2
$ws1.Range("H9").Value = 8.5
$ws1.Range("L9").Value = 1.11

# Row 10 (W18)
$ws1.Range("D10").Value = 2
$ws1.Range("H10").Value = 7.5
$ws1.Range("L10").Value = 1.12

# Row 11 (W19)
$ws1.Range("D11").Value = 2
$ws1.Range("H11").Value = 6.5
$ws1.Range("L11").Value = 1.19

# Row 12 (W20)
$ws1.Range("D12").Value = 2
$ws1.Range("H12").Value = 5.5
$ws1.Range("L12").Value = 1.14

# Row 13 (W21)
$ws1.Range("D13").Value = 2
$ws1.Range("H13").Value = 4.5
$ws1.Range("L13").Value = 0.89

# Row 14 (W22)
$ws1.Range("D14").Value = 2
$ws1.Range("H14").Value = 3.5
$ws1.Range("L14").Value = 1.09

# Row 15 (W23)
$ws1.Range("D15").Value = 2
$ws1.Range("H15").Value = 2.5
$ws1.Range("L15").Value = 1.1

# Row 16 (W24)
$ws1.Range("H16").Value = 3
$ws1.Range("L16").Value = 1.15

# Row 17 (W25)
$ws1.Range("H17").Value = 2
$ws1.Range("L17").Value = 1.08

# --- Summary sheet ----------------------------------------------------------

# Total Forecast (16 Weeks)
$ws2.Range("B9").Value = 28

# Total Forecast (4 Weeks)
$ws2.Range("B11").Value = 6
